$wb = $excel.ActiveWorkbook
$srcSheet = $wb.Worksheets.Item(1)

# Duplicate the existing sheet (preserves columns, styles, formulas,
# external-workbook references, etc.) and place it right after the source.
$srcSheet.Copy($null, $srcSheet)

$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "Transmittals_New_ActionOverDue"

# Update the row-2 values that differ from the "ActionRequired" sheet.
# (Order matters for shared-string append order: Overdue is registered
# before LATFULPP-4.)
$newSheet.Range("P2").Value = "Overdue"
$newSheet.Range("A2").Value = "LATFULPP-4"

# Match the recorded selection on the new (now active) sheet.
$newSheet.Range("A2").Select()

# New sheet is the active tab.
$newSheet.Activate()
